# Tagesbericht und Pflichtenheft aktualisiert
#
# - Row 14 (Pflichtenheft bearbeiten): Fortschritt 10% -> 30%
# - Row 16 (Datenbankanbindung implementieren): Uhrzeit "10.20 Uhr - "
#   -> "10.20 Uhr -  12.00 Uhr" (the session has now ended, so the end time
#   was filled in)
# - The active selection moves to E14, the cell that was last edited

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fortschritt (progress) for "Pflichtenheft bearbeiten" raised from 10% to 30%
$ws.Range("E14").Value = 0.3

# Uhrzeit for "Datenbankanbindung implementieren" now includes the end time
$ws.Range("D16").Value = "10.20 Uhr -  12.00 Uhr"

# Reflect the cell that was last worked on as the active selection
$ws.Range("E14").Select()
